$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("B4").Value = 1708
$ws.Range("B7").Value = 23340
$ws.Range("K4").Value = 1763
$ws.Range("K7").Value = 27554
$ws.Range("L2").Value = 2137
$ws.Range("L3").Value = 2147
$ws.Range("L4").Value = 594
$ws.Range("L6").Value = 1937
$ws.Range("L7").Value = 6938

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B63").Value = 413
$ws.Range("B101").Value = 23340
$ws.Range("J29").Value = 1552
$ws.Range("J63").Value = 217
$ws.Range("K63").Value = 88
$ws.Range("K101").Value = 27554
$ws.Range("L7").Value = 228
$ws.Range("L8").Value = 435
$ws.Range("L11").Value = 127
$ws.Range("L15").Value = 48
$ws.Range("L16").Value = 15
$ws.Range("L19").Value = 196
$ws.Range("L23").Value = 71
$ws.Range("L24").Value = 15
$ws.Range("L29").Value = 361
$ws.Range("L33").Value = 313
$ws.Range("L39").Value = 2
$ws.Range("L42").Value = 216
$ws.Range("L48").Value = 97
$ws.Range("L50").Value = 41
$ws.Range("L51").Value = 80
$ws.Range("L55").Value = 62
$ws.Range("L57").Value = 30
$ws.Range("L60").Value = 37
$ws.Range("L64").Value = 51
$ws.Range("L65").Value = 133
$ws.Range("L74").Value = 7
$ws.Range("L76").Value = 73
$ws.Range("L79").Value = 189
$ws.Range("L83").Value = 163
$ws.Range("L85").Value = 365
$ws.Range("L88").Value = 95
$ws.Range("L89").Value = 91
$ws.Range("L90").Value = 67
$ws.Range("L91").Value = 101
$ws.Range("L94").Value = 81
$ws.Range("L97").Value = 65
$ws.Range("L99").Value = 108
$ws.Range("L101").Value = 6938

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 66
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 43
$ws.Range("L4").Value = 11
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 150
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 365

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 150
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 435

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 313

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J4").Value = 84
$ws.Range("J7").Value = 1552
$ws.Range("L3").Value = 127
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 361

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 64
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 61
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 2

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L2").Value = 11
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 7
